$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 242, pushing existing rows 242:255 down to 243:256
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new data record
$ws.Cells.Item(242, 1).Value = 10
$ws.Cells.Item(242, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(242, 3).Value = "La Araucanía"
$ws.Cells.Item(242, 4).Value = 44585
$ws.Cells.Item(242, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(242, 5).Value = 9
$ws.Cells.Item(242, 6).Value = 100114013
$ws.Cells.Item(242, 7).Value = "Zanahoria"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 125
$ws.Cells.Item(242, 11).Value = 8000
$ws.Cells.Item(242, 12).Value = 8000
$ws.Cells.Item(242, 13).Value = 8000
$ws.Cells.Item(242, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 400
$ws.Cells.Item(242, 17).Value = 20
$ws.Cells.Item(242, 18).Value = "Hortaliza"
